$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 30. This shifts the existing rows 30-44
# (weekly "Ají" price records) down to rows 31-45, growing the used
# range from A1:R44 to A1:R45, and carries the date-number-format from
# row 30 down onto the new blank row (matching column D's style).
$ws.Rows("30:30").Insert()

# Populate the newly inserted row 30 with this week's record.
$ws.Range("A30").Value = 1
$ws.Range("B30").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C30").Value = "Arica y Parinacota"
$ws.Range("D30").Value = 44466
$ws.Range("E30").Value = 15
$ws.Range("F30").Value = 100112021
$ws.Range("G30").Value = "Ají"
$ws.Range("H30").Value = "Inferno"
$ws.Range("I30").Value = "Primera"
$ws.Range("J30").Value = 120
$ws.Range("K30").Value = 36000
$ws.Range("L30").Value = 38000
$ws.Range("M30").Value = 37000
$ws.Range("N30").Value = "$/caja 15 kilos"
$ws.Range("O30").Value = "Región de Arica y Parinacota"
$ws.Range("P30").Value = 2467
$ws.Range("Q30").Value = 15
$ws.Range("R30").Value = "Hortaliza"
